$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("History.html")
$ws.Rows.Item(3).Resize(2).Insert()
Write-Host "Row3 height:" $ws.Rows.Item(3).RowHeight
Write-Host "Row4 height:" $ws.Rows.Item(4).RowHeight
